$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing date cells (rows 13 and 20) from 06/18/2019 to 08/08/2019
$ws.Range("D13:G13").Value2 = 43685
$ws.Range("D20:G20").Value2 = 43685

# Rows 14-19 and 21-23 previously held an "OK" marker / empty cells with a
# different fill style; bring them in line with the completed rows by
# copying the date format from row 13 and filling in the new date value.
$ws.Range("D13").Copy()
$ws.Range("D14:G19").PasteSpecial(-4122)
$ws.Range("D21:G23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D14:G19").Value2 = 43685
$ws.Range("D21:G23").Value2 = 43685

# Re-touch the other merged cell ranges so the workbook's internal merge
# list matches the order Excel produced when it re-saved the file.
$touchOrder = @("A24:A27","A2:A5","A6:A7","A8:A12","A13:A19","A20:A23")
foreach ($ref in $touchOrder) {
  $ws.Range($ref).UnMerge()
  $ws.Range($ref).Merge()
}

# Selection state as left by the editor after finishing this section
$ws.Range("E26").Select() | Out-Null
